$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 859, shifting existing rows 859:900 down to 860:901
$ws.Rows(859).Insert()

# Populate the newly inserted row 859 with the new data point.
# Force column A to text format first so the date-like string isn't
# auto-converted into a date serial number, then restore the default
# "Normal" style so the cell doesn't end up with a stray explicit style
# (matching the plain, unstyled look of the surrounding data rows).
$ws.Range("A859").NumberFormat = "@"
$ws.Range("A859").Value = "2026/02/22"
$ws.Range("A859").Style = "Normal"
$ws.Range("B859").Value = "日"
$ws.Range("C859").Value = 13
$ws.Range("D859").Value = 201
